# Adds row 24 to the "Artfynd" sheet, mirroring the new species-observation
# record introduced by the source diff (dimension A1:AY23 -> A1:AY24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric cells -----------------------------------------------------
$ws.Range("A24").Value()  = 112222749
$ws.Range("B24").Value()  = 89405
$ws.Range("E24").Value()  = 1202
$ws.Range("Q24").Value()  = 558131
$ws.Range("R24").Value()  = 6628068
$ws.Range("S24").Value()  = 25

# --- Plain text cells ----------------------------------------------------
$ws.Range("C24").Value()  = "Ovaliderad"
$ws.Range("D24").Value()  = "NT"
$ws.Range("F24").Value()  = "Ullticka"
$ws.Range("G24").Value()  = "Phellinidium ferrugineofuscum"
$ws.Range("H24").Value()  = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P24").Value()  = "Månses hål S, Vstm"
$ws.Range("T24").Value()  = "Västmanland"
$ws.Range("U24").Value()  = "Surahammar"
$ws.Range("V24").Value()  = "Västmanland"
$ws.Range("W24").Value()  = "Ramnäs"
$ws.Range("AI24").Value() = "Blandskog"
$ws.Range("AO24").Value() = "Gran"
$ws.Range("AW24").Value() = "Tom Sävström"
$ws.Range("AX24").Value() = "Tom Sävström"

# --- Date-like text (kept as literal text, not converted to a date) ------
# A leading apostrophe forces Excel to treat the value as text, which is
# what the source file stores for these "Startdatum"/"Slutdatum" fields.
# The apostrophe prefix also flips on a "quote prefix" cell style, so the
# style is reset back to Normal immediately afterwards to avoid leaving any
# formatting difference behind.
$ws.Range("Y24").Value()  = "'2023-09-19"
$ws.Range("Y24").Style()  = "Normal"
$ws.Range("AA24").Value() = "'2023-09-19"
$ws.Range("AA24").Style() = "Normal"

# --- Boolean cells ---------------------------------------------------------
$ws.Range("AD24").Value() = $false
$ws.Range("AE24").Value() = $false
$ws.Range("AG24").Value() = $false

# --- Present-but-empty text cells ------------------------------------------
# The source row has several columns that exist (t="inlineStr") but carry no
# text. Assigning a bare "" clears/removes a cell in this object model, so a
# leading apostrophe is used again to store a real (empty) text value while
# still reading back as "" — the quote-prefix style is then reset to Normal,
# same as above.
$ws.Range("I24").Value()  = "'"
$ws.Range("I24").Style()  = "Normal"
$ws.Range("J24").Value()  = "'"
$ws.Range("J24").Style()  = "Normal"
$ws.Range("K24").Value()  = "'"
$ws.Range("K24").Style()  = "Normal"
$ws.Range("N24").Value()  = "'"
$ws.Range("N24").Style()  = "Normal"
$ws.Range("AF24").Value() = "'"
$ws.Range("AF24").Style() = "Normal"
$ws.Range("AT24").Value() = "'"
$ws.Range("AT24").Style() = "Normal"
$ws.Range("AY24").Value() = "'"
$ws.Range("AY24").Style() = "Normal"

Write-Output "Row 24 added."
